$wb = $excel.ActiveWorkbook

# --- 1. Rename the "ja-jp" worksheet to "it-it" ---
$wsJp = $wb.Worksheets.Item("ja-jp")
$wsJp.Name = "it-it"

$wsOverview = $wb.Worksheets.Item("Overview")

# --- 2. Update the handoff/handback datetime values on the renamed sheet ---
$wsJp.Range("E2:E5").Value = "2016-03-11 01:04:01"
$wsJp.Range("H2:H5").Value = "2016-03-17 01:40:10"

# --- 3. Rename the table ("ListObject") that lives on the it-it sheet ---
$loItIt = $wsJp.ListObjects.Item(1)
$loItIt.Name = "it-it"

# --- 4. Update the Overview sheet's header text + the Overview table's 2nd column name ---
$wsOverview.Range("B1").Value = "it-it"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Unlist()

$loOverview2 = $wsOverview.ListObjects.Add(1, $wsOverview.Range("A1:C1"), 0, 0)
$loOverview2.ShowHeaders = $false
$loOverview2.Resize($wsOverview.Range("A1:C1"))
$loOverview2.Name = "Overview"
$loOverview2.TableStyle = "TableStyleMedium9"

# Re-adding the table via .Add() stamps a literal header value ("Column3") into the
# previously-empty C1 cell; clear it back out so the sheet data goes back to untouched.
$wsOverview.Range("C1").ClearContents()
